$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8
$ws.Range("I2").Value = 1.4
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.87
$ws.Range("R2").Value = 2.03
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("W2").Value = 19
$ws.Range("Y2").Value = 23
$ws.Range("Z2").Value = 101
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 21
$ws.Range("AG2").Value = 401
$ws.Range("AH2").Value = 6.5
$ws.Range("AK2").Value = 9
$ws.Range("AN2").Value = 8.5
$ws.Range("AS2").Value = 351
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 3.25
$ws.Range("BC2").Value = 151
$ws.Range("G3").Value = 2.6
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3.25
$ws.Range("L3").Value = 3.5
$ws.Range("W3").Value = 7.5
$ws.Range("AA3").Value = 21
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 11
$ws.Range("AK3").Value = 29
$ws.Range("AW3").Value = 4.75
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 5
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 2
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 17
$ws.Range("AN5").Value = 3.6
$ws.Range("AP5").Value = 15
$ws.Range("BA5").Value = 101
$ws.Range("G9").Value = 4.35
$ws.Range("H9").Value = 2.72
$ws.Range("I9").Value = 2.07
$ws.Range("J9").Value = 4.85
$ws.Range("L9").Value = 2.8
$ws.Range("U9").Value = 2.3
$ws.Range("V9").Value = 1.55
$ws.Range("W9").Value = 8.5
$ws.Range("Y9").Value = 15
$ws.Range("AA9").Value = 55
$ws.Range("AB9").Value = 75
$ws.Range("AH9").Value = 4.8
$ws.Range("AI9").Value = 8.25
$ws.Range("AJ9").Value = 9.75
$ws.Range("AK9").Value = 19.5
$ws.Range("AL9").Value = 24
$ws.Range("AO9").Value = 27
$ws.Range("AP9").Value = 37
$ws.Range("AQ9").Value = 175
$ws.Range("AT9").Value = 2.18
$ws.Range("AW9").Value = 3.7
$ws.Range("AY9").Value = 26
$ws.Range("AZ9").Value = 55
$ws.Range("BA9").Value = 120
$ws.Range("BB9").Value = 500
$ws.Range("H10").Value = 2.8
$ws.Range("J10").Value = 3.4
$ws.Range("AC10").Value = 6
$ws.Range("AM10").Value = 41
$ws.Range("BB10").Value = 351
$ws.Range("G13").Value = 1.85
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 4.33
$ws.Range("J13").Value = 2.6
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 8
$ws.Range("O13").Value = 1.4
$ws.Range("P13").Value = 2.75
$ws.Range("Q13").Value = 2.25
$ws.Range("R13").Value = 1.62
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.75
$ws.Range("X13").Value = 8
$ws.Range("AC13").Value = 8
$ws.Range("AE13").Value = 17
$ws.Range("AG13").Value = 451
$ws.Range("AM13").Value = 41
$ws.Range("AN13").Value = 3.75
$ws.Range("AY13").Value = 34
$ws.Range("AZ13").Value = 81
$ws.Range("BB13").Value = 301
$ws.Range("N17").Value = 8.5
$ws.Range("X17").Value = 17
$ws.Range("G20").Value = 1.38
$ws.Range("H20").Value = 4.8
$ws.Range("I20").Value = 7.4
$ws.Range("J20").Value = 1.85
$ws.Range("K20").Value = 2.4
$ws.Range("L20").Value = 6.5
$ws.Range("N20").Value = 9.25
$ws.Range("O20").Value = 1.22
$ws.Range("P20").Value = 3.95
$ws.Range("Q20").Value = 1.7
$ws.Range("R20").Value = 2.12
$ws.Range("S20").Value = 1.35
$ws.Range("T20").Value = 3.15
$ws.Range("U20").Value = 1.98
$ws.Range("Z20").Value = 9.25
$ws.Range("AA20").Value = 12
$ws.Range("AC20").Value = 9.25
$ws.Range("AD20").Value = 9.75
$ws.Range("AE20").Value = 23
$ws.Range("AH20").Value = 17
$ws.Range("AJ20").Value = 25
$ws.Range("AK20").Value = 200
$ws.Range("AO20").Value = 6.2
$ws.Range("AQ20").Value = 17
$ws.Range("AR20").Value = 45
$ws.Range("AU20").Value = 8.75
$ws.Range("AW20").Value = 8.25
$ws.Range("BA20").Value = 250
$ws.Range("G23").Value = 2.4
$ws.Range("I23").Value = 2.75
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = 3.25
$ws.Range("Q23").Value = 1.62
$ws.Range("R23").Value = 2.25
$ws.Range("W23").Value = 11
$ws.Range("Y23").Value = 9.5
$ws.Range("Z23").Value = 23
$ws.Range("AB23").Value = 23
$ws.Range("AH23").Value = 12
$ws.Range("AI23").Value = 15
$ws.Range("AK23").Value = 29
$ws.Range("AM23").Value = 23
$ws.Range("AN23").Value = 4.75
$ws.Range("AO23").Value = 13
$ws.Range("AZ23").Value = 41

Write-Output "Applied 139 cell updates"
